$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("begroting")

# Bug fix: booked hours formula needs to include two more entries (34 + 16)
$ws.Range("B22").Formula = "=8+34+34+16"

# Update the active selection to reflect where the user ended up after the edit
$ws.Range("D20").Select()
